$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.136.21"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "3.614.44"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.04"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "193.69"
$ws.Range("E6").Value = "  +2.50%  "
$ws.Range("D7").Value = "3.609.75"
$ws.Range("E7").Value = "  -1.09%  "
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.152"
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.43"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000290"
$ws.Range("E13").Value = "  +8.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.04"
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("D15").Value = "4.183.46"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").Value = "3.618.78"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "67.989.87"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.54"
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("E21").Value = "  -1.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "404.89"
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.24"
$ws.Range("E23").Value = "  +22.57%  "
$ws.Range("E24").Value = "  -3.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.05"
$ws.Range("E25").Value = "  -1.49%  "
$ws.Range("E26").Value = "  +1.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.97"
$ws.Range("E27").Value = "  +8.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.58"
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.12"
$ws.Range("E30").Value = "  +14.48%  "
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "688.63"
$ws.Range("E33").Value = "  +13.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.27"
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("E35").Value = "  +2.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "64.70"
$ws.Range("E36").Value = "  -4.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.72"
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.422"
$ws.Range("E38").Value = "  +8.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").Value = "0.0₃0792"
$ws.Range("E40").Value = "  +4.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.94"
$ws.Range("E41").Value = "  +17.98%  "
$ws.Range("E42").Value = "  +9.66%  "
$ws.Range("D43").Value = "3.177.79"
$ws.Range("E43").Value = "  +15.80%  "
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("E46").Value = "  +0.77%  "
$ws.Range("E47").Value = "  -1.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.84"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.11"
$ws.Range("E49").Value = "  -3.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "143.35"
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.58"
$ws.Range("E51").Value = "  +2.91%  "
